# Apply the per-row "Price" (D) and "Volume(1h)" (E) updates from the
# coinranking crypto-list refresh. Values are written as literal text so
# they match the original inline-string cell contents (e.g. thousands-dot
# formatted prices like "29.550.55" and padded percentages like "  +4.29%  ").
#
# A handful of Price values look like plain decimal numbers (e.g. "0.994").
# Those are written with a leading apostrophe, Excel's standard "force text"
# quote-prefix, so they don't get auto-converted to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.550.55"
$ws.Range("E2").Value = "  +4.29%  "
$ws.Range("D3").Value = "1.598.22"
$ws.Range("E3").Value = "  +3.12%  "
$ws.Range("D4").Value = "'0.994"
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "'213.30"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "'0.513"
$ws.Range("E6").Value = "  +6.88%  "
$ws.Range("D7").Value = "'0.993"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").Value = "'26.72"
$ws.Range("E8").Value = "  +12.17%  "
$ws.Range("E9").Value = "  +3.31%  "
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("D12").Value = "1.826.99"
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("D13").Value = "1.591.10"
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("D14").Value = "29.568.40"
$ws.Range("E14").Value = "  +4.40%  "
$ws.Range("D15").Value = "'0.528"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").Value = "'3.74"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").Value = "'63.31"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").Value = "'242.23"
$ws.Range("E18").Value = "  +6.52%  "
$ws.Range("D19").Value = "'7.57"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").Value = "0.0₃0693"
$ws.Range("D21").Value = "'0.993"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").Value = "'4.04"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("E23").Value = "  +4.27%  "
$ws.Range("D24").Value = "'2.10"
$ws.Range("E24").Value = "  +3.74%  "
$ws.Range("D25").Value = "'155.01"
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("D26").Value = "'15.27"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("E27").Value = "  +5.38%  "
$ws.Range("D28").Value = "'6.39"
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("D29").Value = "'0.994"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "'0.0472"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").Value = "1.431.49"
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("D35").Value = "'1.04"
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("E37").Value = "  +9.14%  "
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").Value = "'0.535"
$ws.Range("E40").Value = "  +4.93%  "
$ws.Range("D41").Value = "'1.96"
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D42").Value = "'54.08"
$ws.Range("E42").Value = "  +28.66%  "
$ws.Range("D43").Value = "'0.802"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("D44").Value = "'0.992"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").Value = "'65.53"
$ws.Range("E46").Value = "  +5.82%  "
$ws.Range("D47").Value = "'5.36"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "1.737.72"
$ws.Range("E48").Value = "  +3.18%  "
$ws.Range("D49").Value = "'86.24"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "'0.837"
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("E51").Value = "  +1.72%  "
